$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 54, shifting existing rows 54..128 down to 55..129
$ws.Rows.Item(54).Insert()

$ws.Cells.Item(54, 1).Value = 11
$ws.Cells.Item(54, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(54, 3).Value = "Bíobío"
$ws.Cells.Item(54, 4).Value = 44721
$ws.Cells.Item(54, 5).Value = 8
$ws.Cells.Item(54, 6).Value = 100112032
$ws.Cells.Item(54, 7).Value = "Zapallo italiano"
$ws.Cells.Item(54, 8).Value = "Huracán"
$ws.Cells.Item(54, 9).Value = "Primera"
$ws.Cells.Item(54, 10).Value = 180
$ws.Cells.Item(54, 11).Value = 7500
$ws.Cells.Item(54, 12).Value = 8000
$ws.Cells.Item(54, 13).Value = 7778
$ws.Cells.Item(54, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(54, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(54, 16).Value = 130
$ws.Cells.Item(54, 17).Value = 60
$ws.Cells.Item(54, 18).Value = "Hortaliza"
